# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across several
# job sheets in the Seraph_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value2 = 1067.1555  # ALC!H15 was 1018.02325
$ws.Cells.Item(15, 9).Value2 = 1067.1555  # ALC!I15 was 1018.02325
$ws.Cells.Item(15, 11).Value2 = 3201.4665  # ALC!K15 was 3054.06975
$ws.Cells.Item(15, 13).Value2 = -3032.4665  # ALC!M15 was -2885.06975

$ws.Cells.Item(17, 8).Value2 = 1916.5  # ALC!H17 was 1732.6364
$ws.Cells.Item(17, 9).Value2 = 0  # ALC!I17 was 1063
$ws.Cells.Item(17, 10).Value2 = 1916.5  # ALC!J17 was 1799.6
$ws.Cells.Item(17, 11).Value2 = 0  # ALC!K17 was 3189
$ws.Cells.Item(17, 12).Value2 = 5749.5  # ALC!L17 was 5398.799999999999
$ws.Cells.Item(17, 13).ClearContents()  # ALC!M17 was -3021
$ws.Cells.Item(17, 14).Value2 = -6085.5  # ALC!N17 was -5734.799999999999

$ws.Cells.Item(18, 8).Value2 = 1814.3334  # ALC!H18 was 2124.5
$ws.Cells.Item(18, 9).Value2 = 1814.3334  # ALC!I18 was 2124.5
$ws.Cells.Item(18, 11).Value2 = 1814.3334  # ALC!K18 was 2124.5
$ws.Cells.Item(18, 13).Value2 = -1530.3334  # ALC!M18 was -1840.5

$ws.Cells.Item(28, 8).Value2 = 371.54544  # ALC!H28 was 401.8
$ws.Cells.Item(28, 9).Value2 = 371.54544  # ALC!I28 was 401.8
$ws.Cells.Item(28, 11).Value2 = 371.54544  # ALC!K28 was 401.8
$ws.Cells.Item(28, 13).Value2 = 113.45456  # ALC!M28 was 83.19999999999999

$ws.Cells.Item(32, 8).Value2 = 1803.9  # ALC!H32 was 1959.6666
$ws.Cells.Item(32, 10).Value2 = 3036  # ALC!J32 was 3914
$ws.Cells.Item(32, 12).Value2 = 3036  # ALC!L32 was 3914
$ws.Cells.Item(32, 14).Value2 = -3688  # ALC!N32 was -4566

$ws.Cells.Item(43, 8).Value2 = 12464.353  # ALC!H43 was 12309.5
$ws.Cells.Item(43, 9).Value2 = 10808.637  # ALC!I43 was 10387.25
$ws.Cells.Item(43, 10).Value2 = 15499.833  # ALC!J43 was 19998.5
$ws.Cells.Item(43, 11).Value2 = 10808.637  # ALC!K43 was 10387.25
$ws.Cells.Item(43, 12).Value2 = 15499.833  # ALC!L43 was 19998.5
$ws.Cells.Item(43, 13).Value2 = -10739.637  # ALC!M43 was -10318.25
$ws.Cells.Item(43, 14).Value2 = -15637.833  # ALC!N43 was -20136.5

$ws.Cells.Item(99, 8).Value2 = 302.42856  # ALC!H99 was 341.6
$ws.Cells.Item(99, 9).Value2 = 287.6  # ALC!I99 was 302
$ws.Cells.Item(99, 10).Value2 = 339.5  # ALC!J99 was 500
$ws.Cells.Item(99, 11).Value2 = 862.8000000000001  # ALC!K99 was 906
$ws.Cells.Item(99, 12).Value2 = 1018.5  # ALC!L99 was 1500
$ws.Cells.Item(99, 13).Value2 = 635.1999999999999  # ALC!M99 was 592
$ws.Cells.Item(99, 14).Value2 = -4014.5  # ALC!N99 was -4496

$ws.Cells.Item(135, 8).Value2 = 1149.88  # ALC!H135 was 1150.28
$ws.Cells.Item(135, 9).Value2 = 645.1579  # ALC!I135 was 645.6842
$ws.Cells.Item(135, 11).Value2 = 5806.4211  # ALC!K135 was 5811.1578
$ws.Cells.Item(135, 13).Value2 = -3271.4211  # ALC!M135 was -3276.1578

$ws.Cells.Item(137, 8).Value2 = 1055.375  # ALC!H137 was 1058.5625
$ws.Cells.Item(137, 9).Value2 = 925.5  # ALC!I137 was 929.75
$ws.Cells.Item(137, 11).Value2 = 2776.5  # ALC!K137 was 2789.25
$ws.Cells.Item(137, 13).Value2 = -226.5  # ALC!M137 was -239.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value2 = 645.2258  # ARM!H74 was 662.5862
$ws.Cells.Item(74, 9).Value2 = 662.2414  # ARM!I74 was 682.14813
$ws.Cells.Item(74, 11).Value2 = 662.2414  # ARM!K74 was 682.14813
$ws.Cells.Item(74, 13).Value2 = 211.7586  # ARM!M74 was 191.85187

$ws.Cells.Item(77, 8).Value2 = 645.2258  # ARM!H77 was 662.5862
$ws.Cells.Item(77, 9).Value2 = 662.2414  # ARM!I77 was 682.14813
$ws.Cells.Item(77, 11).Value2 = 3311.207  # ARM!K77 was 3410.74065
$ws.Cells.Item(77, 13).Value2 = 1056.793  # ARM!M77 was 957.2593499999998

$ws.Cells.Item(132, 8).Value2 = 14586.75  # ARM!H132 was 15426
$ws.Cells.Item(132, 9).Value2 = 2384.8572  # ARM!I132 was 2414.6155
$ws.Cells.Item(132, 11).Value2 = 7154.571599999999  # ARM!K132 was 7243.8465
$ws.Cells.Item(132, 13).Value2 = -4624.571599999999  # ARM!M132 was -4713.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value2 = 111280  # BSM!H140 was 107520
$ws.Cells.Item(140, 10).Value2 = 111280  # BSM!J140 was 107520
$ws.Cells.Item(140, 12).Value2 = 111280  # BSM!L140 was 107520
$ws.Cells.Item(140, 14).Value2 = -121640  # BSM!N140 was -117880

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value2 = 1894.1052  # CRP!H58 was 1925.75
$ws.Cells.Item(58, 9).Value2 = 1118.2142  # CRP!I58 was 1102.3462
$ws.Cells.Item(58, 11).Value2 = 1118.2142  # CRP!K58 was 1102.3462
$ws.Cells.Item(58, 13).Value2 = -915.2141999999999  # CRP!M58 was -899.3462

$ws.Cells.Item(107, 8).Value2 = 980.9231  # CRP!H107 was 867.46155
$ws.Cells.Item(107, 9).Value2 = 603.5  # CRP!I107 was 594.7778
$ws.Cells.Item(107, 10).Value2 = 1584.8  # CRP!J107 was 1481
$ws.Cells.Item(107, 11).Value2 = 603.5  # CRP!K107 was 594.7778
$ws.Cells.Item(107, 12).Value2 = 1584.8  # CRP!L107 was 1481
$ws.Cells.Item(107, 13).Value2 = 1316.5  # CRP!M107 was 1325.2222
$ws.Cells.Item(107, 14).Value2 = -5424.8  # CRP!N107 was -5321

$ws.Cells.Item(134, 8).Value2 = 2400.75  # CRP!H134 was 2465.75
$ws.Cells.Item(134, 9).Value2 = 2239.3428  # CRP!I134 was 2328.4856
$ws.Cells.Item(134, 11).Value2 = 6718.028399999999  # CRP!K134 was 6985.4568
$ws.Cells.Item(134, 13).Value2 = -4183.028399999999  # CRP!M134 was -4450.4568

$ws.Cells.Item(136, 8).Value2 = 1894.1052  # CRP!H136 was 1925.75
$ws.Cells.Item(136, 9).Value2 = 1118.2142  # CRP!I136 was 1102.3462
$ws.Cells.Item(136, 11).Value2 = 3354.6426  # CRP!K136 was 3307.0386
$ws.Cells.Item(136, 13).Value2 = -804.6425999999997  # CRP!M136 was -757.0385999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(134, 8).Value2 = 3805.4  # CUL!H134 was 4004.2
$ws.Cells.Item(134, 9).Value2 = 2009.3334  # CUL!I134 was 2014.5
$ws.Cells.Item(134, 10).Value2 = 6499.5  # CUL!J134 was 5330.6665
$ws.Cells.Item(134, 11).Value2 = 6028.0002  # CUL!K134 was 6043.5
$ws.Cells.Item(134, 12).Value2 = 19498.5  # CUL!L134 was 15991.9995
$ws.Cells.Item(134, 13).Value2 = -958.0002000000004  # CUL!M134 was -973.5
$ws.Cells.Item(134, 14).Value2 = -29638.5  # CUL!N134 was -26131.9995

$ws.Cells.Item(138, 8).Value2 = 7000  # CUL!H138 was 0
$ws.Cells.Item(138, 10).Value2 = 7000  # CUL!J138 was 0
$ws.Cells.Item(138, 12).Value2 = 21000  # CUL!L138 was 0
$ws.Cells.Item(138, 14).Value2 = -31280  # CUL!N138 was None

$ws.Cells.Item(139, 8).Value2 = 2951.6  # CUL!H139 was 5006.25
$ws.Cells.Item(139, 9).Value2 = 1749.5  # CUL!I139 was 2250
$ws.Cells.Item(139, 10).Value2 = 7760  # CUL!J139 was 7762.5
$ws.Cells.Item(139, 11).Value2 = 5248.5  # CUL!K139 was 6750
$ws.Cells.Item(139, 12).Value2 = 23280  # CUL!L139 was 23287.5
$ws.Cells.Item(139, 13).Value2 = -108.5  # CUL!M139 was -1610
$ws.Cells.Item(139, 14).Value2 = -33560  # CUL!N139 was -33567.5

$ws.Cells.Item(140, 8).Value2 = 4180  # CUL!H140 was 4802.2
$ws.Cells.Item(140, 9).Value2 = 2869.6  # CUL!I140 was 4015
$ws.Cells.Item(140, 11).Value2 = 8608.799999999999  # CUL!K140 was 12045
$ws.Cells.Item(140, 13).Value2 = -3428.799999999999  # CUL!M140 was -6865

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value2 = 418.58823  # GSM!H2 was 438.5625
$ws.Cells.Item(2, 10).Value2 = 1128.8  # GSM!J2 was 1386.25
$ws.Cells.Item(2, 12).Value2 = 1128.8  # GSM!L2 was 1386.25
$ws.Cells.Item(2, 14).Value2 = -1354.8  # GSM!N2 was -1612.25

$ws.Cells.Item(70, 8).Value2 = 6155.091  # GSM!H70 was 5538.077
$ws.Cells.Item(70, 9).Value2 = 4741.4  # GSM!I70 was 3999.4285
$ws.Cells.Item(70, 11).Value2 = 4741.4  # GSM!K70 was 3999.4285
$ws.Cells.Item(70, 13).Value2 = -4471.4  # GSM!M70 was -3729.4285

$ws.Cells.Item(73, 8).Value2 = 6155.091  # GSM!H73 was 5538.077
$ws.Cells.Item(73, 9).Value2 = 4741.4  # GSM!I73 was 3999.4285
$ws.Cells.Item(73, 11).Value2 = 4741.4  # GSM!K73 was 3999.4285
$ws.Cells.Item(73, 13).Value2 = -3805.4  # GSM!M73 was -3063.4285

$ws.Cells.Item(122, 8).Value2 = 92729.91  # GSM!H122 was 102509.9
$ws.Cells.Item(122, 9).Value2 = 1253.875  # GSM!I122 was 1355.3334
$ws.Cells.Item(122, 10).Value2 = 336666  # GSM!J122 was 254241.75
$ws.Cells.Item(122, 11).Value2 = 3761.625  # GSM!K122 was 4066.0002
$ws.Cells.Item(122, 12).Value2 = 1009998  # GSM!L122 was 762725.25
$ws.Cells.Item(122, 13).Value2 = -1311.625  # GSM!M122 was -1616.0002
$ws.Cells.Item(122, 14).Value2 = -1014898  # GSM!N122 was -767625.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value2 = 0  # LTW!H36 was 81999
$ws.Cells.Item(36, 10).Value2 = 0  # LTW!J36 was 81999
$ws.Cells.Item(36, 12).Value2 = 0  # LTW!L36 was 81999
$ws.Cells.Item(36, 14).ClearContents()  # LTW!N36 was -83123

$ws.Cells.Item(136, 8).Value2 = 7357.615  # LTW!H136 was 7364.615
$ws.Cells.Item(136, 9).Value2 = 7296.5454  # LTW!I136 was 7304.8184
$ws.Cells.Item(136, 11).Value2 = 21889.6362  # LTW!K136 was 21914.4552
$ws.Cells.Item(136, 13).Value2 = -19339.6362  # LTW!M136 was -19364.4552

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value2 = 500  # WVR!H2 was 437.5
$ws.Cells.Item(2, 9).Value2 = 0  # WVR!I2 was 250
$ws.Cells.Item(2, 11).Value2 = 0  # WVR!K2 was 250
$ws.Cells.Item(2, 13).ClearContents()  # WVR!M2 was -138

$ws.Cells.Item(62, 8).Value2 = 7478.087  # WVR!H62 was 7558.909
$ws.Cells.Item(62, 10).Value2 = 7685.05  # WVR!J62 was 7789.5264
$ws.Cells.Item(62, 12).Value2 = 7685.05  # WVR!L62 was 7789.5264
$ws.Cells.Item(62, 14).Value2 = -8933.049999999999  # WVR!N62 was -9037.526399999999

$ws.Cells.Item(65, 8).Value2 = 7478.087  # WVR!H65 was 7558.909
$ws.Cells.Item(65, 10).Value2 = 7685.05  # WVR!J65 was 7789.5264
$ws.Cells.Item(65, 12).Value2 = 38425.25  # WVR!L65 was 38947.632
$ws.Cells.Item(65, 14).Value2 = -44665.25  # WVR!N65 was -45187.632

$ws.Cells.Item(70, 8).Value2 = 53000  # WVR!H70 was 50500
$ws.Cells.Item(70, 9).Value2 = 0  # WVR!I70 was 10000
$ws.Cells.Item(70, 10).Value2 = 53000  # WVR!J70 was 58600
$ws.Cells.Item(70, 11).Value2 = 0  # WVR!K70 was 10000
$ws.Cells.Item(70, 12).Value2 = 53000  # WVR!L70 was 58600
$ws.Cells.Item(70, 13).ClearContents()  # WVR!M70 was -9685
$ws.Cells.Item(70, 14).Value2 = -53630  # WVR!N70 was -59230

$ws.Cells.Item(73, 8).Value2 = 53000  # WVR!H73 was 50500
$ws.Cells.Item(73, 9).Value2 = 0  # WVR!I73 was 10000
$ws.Cells.Item(73, 10).Value2 = 53000  # WVR!J73 was 58600
$ws.Cells.Item(73, 11).Value2 = 0  # WVR!K73 was 10000
$ws.Cells.Item(73, 12).Value2 = 53000  # WVR!L73 was 58600
$ws.Cells.Item(73, 13).ClearContents()  # WVR!M73 was -8908
$ws.Cells.Item(73, 14).Value2 = -55184  # WVR!N73 was -60784

$ws.Cells.Item(100, 8).Value2 = 4325.25  # WVR!H100 was 4434.3335
$ws.Cells.Item(100, 9).Value2 = 4432.6665  # WVR!I100 was 4650
$ws.Cells.Item(100, 11).Value2 = 8865.333000000001  # WVR!K100 was 9300
$ws.Cells.Item(100, 13).Value2 = -8324.333000000001  # WVR!M100 was -8759

$ws.Cells.Item(107, 8).Value2 = 200.33333  # WVR!H107 was 200
$ws.Cells.Item(107, 9).Value2 = 200.33333  # WVR!I107 was 200
$ws.Cells.Item(107, 11).Value2 = 600.99999  # WVR!K107 was 600
$ws.Cells.Item(107, 13).Value2 = 1319.00001  # WVR!M107 was 1320

$ws.Cells.Item(113, 8).Value2 = 962.5  # WVR!H113 was 686.8
$ws.Cells.Item(113, 9).Value2 = 1000  # WVR!I113 was 585.3
$ws.Cells.Item(113, 10).Value2 = 925  # WVR!J113 was 889.8
$ws.Cells.Item(113, 11).Value2 = 3000  # WVR!K113 was 1755.9
$ws.Cells.Item(113, 12).Value2 = 2775  # WVR!L113 was 2669.4
$ws.Cells.Item(113, 13).Value2 = -830  # WVR!M113 was 414.1000000000001
$ws.Cells.Item(113, 14).Value2 = -7115  # WVR!N113 was -7009.4
